$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #2")

# --- Row 19 ---
$ws.Range("A19").Value2 = 43201
$ws.Range("A19").NumberFormat = "mm-dd-yy"
$ws.Range("A19").HorizontalAlignment = -4131
$ws.Range("B19").Value2 = "Essayer de regler le probleme de la boule qui depasse et qui retourne sur la limite"
$ws.Range("C19").Value2 = "2h30"

# --- Row 20 ---
$ws.Range("A20").Value2 = 43206
$ws.Range("A20").NumberFormat = "mm-dd-yy"
$ws.Range("A20").HorizontalAlignment = -4131
$ws.Range("C20").Value2 = "3h30"
$ws.Range("B20").Value2 = "Essayer de regler le probleme de la boule qui depasse et qui depasse pas"

# --- Row 21 ---
$ws.Range("B21").Value2 = "Regardez des exemples sur le onaccuracychanged"

# --- Row 22 ---
$ws.Range("A22").Value2 = 43208
$ws.Range("A22").NumberFormat = "mm-dd-yy"
$ws.Range("A22").HorizontalAlignment = -4131
$ws.Range("B22").Value2 = "Essayer de faire un score et regardez sur internet pour convertir les données de "
$ws.Range("C22").Value2 = "3h"

# --- Row 23 ---
$ws.Range("B23").Value2 = "l'accéléromètre"

# --- Selection / active cell update ---
$ws.Activate()
$ws.Range("B23").Select()
